$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K15").Value = 0.2305062539156956
$ws.Range("J16").Value = 0.2193215401759246
$ws.Range("I17").Value = 0.2109873117084238
$ws.Range("H18").Value = 0.2077622620068982
$ws.Range("G19").Value = 0.1868984584576193
$ws.Range("F20").Value = 0.2101374940836094
$ws.Range("E21").Value = 0.2201756597651073
$ws.Range("D22").Value = 0.1085991175498651
$ws.Range("C23").Value = 0.130019622424466
$ws.Range("B24").Value = 0.3662627537369125
